$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Buy Value in GBP"
$ws.Range("G1").Value = "Sell Value in GBP"
$ws.Range("J1").Value = "Fee Value in GBP"

$white = 16777215
for ($col = 1; $col -le 13; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    for ($edge = 7; $edge -le 10; $edge++) {
        $cell.Borders.Item($edge).LineStyle = 1
        $cell.Borders.Item($edge).Color = $white
        $cell.Borders.Item($edge).Weight = 2
    }
    Write-Output ("border set for column " + $col)
}
